$d = $word.ActiveDocument
$r = $d.Range(2537, 2544)
Write-Output $r.WordOpenXML
